$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing rows 1-2 (values swap due to removal of "Day 6 (31/05/2019)" entry)
$ws.Range("B1").Value = "Learning using Sass Variables and Nesting"
$ws.Range("B2").Value = "Learning using Sass Mixins Extends and Functions"

# Row 3 stays the same (Day 3 / install scss note)

# Row 4: Day 4 (26/06/2019) / Converting CSS to Sass text - unchanged text, but shared string index shifts only
$ws.Range("A4").Value = "Day 4 (26/06/2019)"
$ws.Range("B4").Value = "Learning Converting Our CSS Code to Sass Variables and Nesting => Implement to the Natours Project"

# Row 5: Day 5 (27/06/2019)
$ws.Range("A5").Value = "Day 5 (27/06/2019)"

# Row 6: new Day 6 (29/06/2019) entry replacing old Day 6 (31/05/2019) slot
$ws.Range("A6").Value = "Day 6 (29/06/2019)"

# Row 5 col B: 7-1 CSS Architecture text (leading space removed) - set after A6 so shared-string order matches
$ws.Range("B5").Value = "Learning how to implementing the 7-1 CSS Architecture with Sass, how to organizing scss file structure"

$ws.Range("B6").Value = "Learning how to building a custom grid with floats, using :not() pseudo"

# Row 7: brand new row - Day 7 (30/06/2019)
$ws.Range("A7").Value = "Day 7 (30/06/2019)"
$ws.Range("B7").Value = "Completing how to building about section with SASS_Natours Project"

# Apply the same style as other A-column cells (center/center) to the new A7 cell
$ws.Range("A7").HorizontalAlignment = -4108
$ws.Range("A7").VerticalAlignment = -4108

# Widen column B to fit new content (target stored width 101.7109375 chars;
# the COM width setter snaps to a 1/6-character grid, so 100.75 is the
# closest input that lands on the nearest reachable stored value)
$ws.Columns("B").ColumnWidth = 100.75
